$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1170
$ws.Range("J2").Value = 3174.5
$ws.Range("L2").Value = 3174.5
$ws.Range("N2").Value = -3400.5
$ws.Range("H12").Value = 723.875
$ws.Range("I12").Value = 1328.75
$ws.Range("J12").Value = 119
$ws.Range("K12").Value = 1328.75
$ws.Range("L12").Value = 119
$ws.Range("M12").Value = -1158.75
$ws.Range("N12").Value = -459
$ws.Range("H41").Value = 322
$ws.Range("I41").Value = 296.57144
$ws.Range("K41").Value = 296.57144
$ws.Range("M41").Value = 143.42856
$ws.Range("H48").Value = 7140
$ws.Range("I48").Value = 1420
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 4260
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -3968
$ws.Range("N48").Value = -30584
$ws.Range("H56").Value = 7140
$ws.Range("I56").Value = 1420
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 4260
$ws.Range("L56").Value = 30000
$ws.Range("M56").Value = -3726
$ws.Range("N56").Value = -31068
$ws.Range("H100").Value = 2456.8
$ws.Range("I100").Value = 2456.8
$ws.Range("K100").Value = 2456.8
$ws.Range("M100").Value = -1915.8
$ws.Range("H107").Value = 362.7
$ws.Range("I107").Value = 362.7
$ws.Range("K107").Value = 362.7
$ws.Range("M107").Value = 1557.3
$ws.Range("H113").Value = 9125
$ws.Range("I113").Value = 8905.625
$ws.Range("J113").Value = 10002.5
$ws.Range("K113").Value = 8905.625
$ws.Range("L113").Value = 10002.5
$ws.Range("M113").Value = -5651.625
$ws.Range("N113").Value = -16510.5
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = $null
$ws.Range("H137").Value = 821.9375
$ws.Range("I137").Value = 650.0909
$ws.Range("J137").Value = 1200
$ws.Range("K137").Value = 1950.2727
$ws.Range("L137").Value = 3600
$ws.Range("M137").Value = 599.7273
$ws.Range("N137").Value = -8700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2132.3872
$ws.Range("I32").Value = 1936.5862
$ws.Range("K32").Value = 1936.5862
$ws.Range("M32").Value = -1649.5862
$ws.Range("H92").Value = 46147.332
$ws.Range("J92").Value = 46147.332
$ws.Range("L92").Value = 46147.332
$ws.Range("N92").Value = -51139.332
$ws.Range("H132").Value = 1198.5
$ws.Range("I132").Value = 1286
$ws.Range("K132").Value = 3858
$ws.Range("M132").Value = -1328
$ws.Range("H134").Value = 79970
$ws.Range("J134").Value = 79970
$ws.Range("L134").Value = 79970
$ws.Range("N134").Value = -90110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 99999
$ws.Range("J51").Value = 99999
$ws.Range("L51").Value = 99999
$ws.Range("N51").Value = -100981
$ws.Range("H105").Value = 2807
$ws.Range("I105").Value = 2941.5
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2941.5
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -1194.5
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 1366.6666
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 820
$ws.Range("N107").Value = -5740
$ws.Range("H134").Value = 1610.8823
$ws.Range("I134").Value = 1562.3334
$ws.Range("K134").Value = 4687.0002
$ws.Range("M134").Value = -2152.0002
$ws.Range("H135").Value = 32666.666
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 32666.666
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 32666.666
$ws.Range("M135").Value = $null
$ws.Range("N135").Value = -42806.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 465.58334
$ws.Range("I22").Value = 485.875
$ws.Range("J22").Value = 425
$ws.Range("K22").Value = 485.875
$ws.Range("L22").Value = 425
$ws.Range("M22").Value = -135.875
$ws.Range("N22").Value = -1125
$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = $null
$ws.Range("H99").Value = 2968
$ws.Range("J99").Value = 1820
$ws.Range("L99").Value = 1820
$ws.Range("N99").Value = -4816
$ws.Range("H107").Value = 356.7143
$ws.Range("I107").Value = 331.55554
$ws.Range("J107").Value = 402
$ws.Range("K107").Value = 331.55554
$ws.Range("L107").Value = 402
$ws.Range("M107").Value = 1588.44446
$ws.Range("N107").Value = -4242
$ws.Range("H126").Value = 2968
$ws.Range("J126").Value = 1820
$ws.Range("L126").Value = 5460
$ws.Range("N126").Value = -10400
$ws.Range("H134").Value = 1622.8
$ws.Range("I134").Value = 1275
$ws.Range("K134").Value = 3825
$ws.Range("M134").Value = -1290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("H14").Value = 2624.25
$ws.Range("I14").Value = 2624.25
$ws.Range("K14").Value = 7872.75
$ws.Range("M14").Value = -7699.75
$ws.Range("H109").Value = 3999
$ws.Range("J109").Value = 3999
$ws.Range("L109").Value = 11997
$ws.Range("N109").Value = -14077
$ws.Range("H140").Value = 920.5
$ws.Range("I140").Value = 800.55554
$ws.Range("K140").Value = 2401.66662
$ws.Range("M140").Value = 2778.33338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5660
$ws.Range("I126").Value = 5433.3335
$ws.Range("K126").Value = 16300.0005
$ws.Range("M126").Value = -13830.0005
$ws.Range("H132").Value = 4253.091
$ws.Range("I132").Value = 4253.091
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12759.273
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10229.273
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H107").Value = 988
$ws.Range("I107").Value = 1011.8571
$ws.Range("K107").Value = 3035.5713
$ws.Range("M107").Value = -1115.5713
